# Replace GroupTC-HS ("grouptc-hash-v2") data with GroupTC-cuckoo values
# Updates column E (raw GroupTC-HS timings) and column I (GroupTC-HS_speedup = Polak / GroupTC-HS)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newE = @{
    2 = 0.9129999999999999
    3 = 2.41
    4 = 7.841000000000001
    5 = 27.214
    6 = 90.235
    7 = 287.543
    8 = 876.48
    9 = 2625.923
}

$newI = @{
    2 = 1.83132530120482
    3 = 2.596265560165976
    4 = 3.019512817242698
    5 = 3.206952303961196
    6 = 3.339469163849948
    7 = 3.438570926783124
    8 = 3.515872581234027
    9 = 3.44856570432568
}

foreach ($row in 2..9) {
    $ws.Cells.Item($row, 5).Value = $newE[$row]
    $ws.Cells.Item($row, 9).Value = $newI[$row]
}
